$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.842.77'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.95%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.868.08'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9998'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '246.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.02%  '
$ws.Range('E6').Value = '  +2.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.000'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.07782'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.15%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3088'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.43%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '23.92'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07851'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('E12').Value = '  +2.91%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.864.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.57%  '
$ws.Range('B14').Value = 'Litecoin'
$ws.Range('C14').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '92.90'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.76%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6973'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.644'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '29.832.43'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008419'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '244.27'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.46%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.114.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('E21').Value = '  +1.45%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.000'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.07%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.667'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.02%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1515'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.26%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.993'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.19'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.46'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.545'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.01%  '
$ws.Range('E30').Value = '  +2.13%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.246'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.98%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.203'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05106'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7911'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.935'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.42%  '
$ws.Range('E36').Value = '  +1.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.711'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.338.42'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.69%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01889'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.753'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9730'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +6.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '6.058'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +11.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '107.18'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9999'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000126'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.77%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '9.833'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.011.98'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '65.51'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.22%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.801'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.15%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.5198'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.44%  '
$ws.Range('E51').Value = '  +2.02%  '
